$wb = $excel.ActiveWorkbook
$wsDatos = $wb.Worksheets.Item("Datos")

# New row for year 2020 on the "Datos" sheet, following the same pattern
# as the other year rows (A: year stored as text, B: numeric indicator value).
# A leading apostrophe forces the numeric-looking "2020" to be kept as text,
# just like every other year label already in column A.
$wsDatos.Range("A22").Value = "'2020"
$wsDatos.Range("B22").Value = 6.2
